$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.63572591109317322
$ws.Range("C2").Value = 0.59358603128570442

$ws.Range("B3").Value = 2.6783732504971551
$ws.Range("C3").Value = 2.9989071078059713

$ws.Range("B4").Value = 1.3921530470386705
$ws.Range("C4").Value = 1.6419451794316096

$ws.Range("B5").Value = 0.44525456386918649
$ws.Range("C5").Value = 0.32797522362799547

$ws.Range("C6").Value = 0.10241660106831009

$ws.Range("B9").Value = 52.774295247091352
$ws.Range("C9").Value = 1.7852239589195316

$ws.Range("C10").Value = 0.763852963893773

$ws.Range("B12").Value = 7.0696844086353581
$ws.Range("C12").Value = 7.1894752448428525

$ws.Range("B15").Value = 0.39532975386303709
$ws.Range("C15").Value = 0.40237362005736332

$ws.Range("B16").Value = 2.0184628303544039
$ws.Range("C16").Value = 2.0633518837148306
